# WIP: added more write surround and the automatic closing
#
# - A2T sheet: row 6 (A10 -> T10 link) gets corrected to A2 -> T3.
# - Components sheet: five new "Agent" rows (A5..A9) are inserted between
#   the existing A4 and A10 rows, pushing the Tool rows down.
# - The active sheet/selection moves from Components to A2T.

$wb = $excel.ActiveWorkbook

# --- A2T: fix the mis-mapped A10/T10 pairing to A2/T3 -----------------
$ws2 = $wb.Worksheets.Item("A2T")
$ws2.Range("B6").Value = "A2"
$ws2.Range("C6").Value = "T3"

# --- Components: insert 5 new Agent rows (A5..A9) before the A10 row --
$ws3 = $wb.Worksheets.Item("Components")
$ws3.Rows("11:15").Insert()

$newAgents = @("A5", "A6", "A7", "A8", "A9")
$newProgress = @(0, 1, 0, 0, 0)

for ($i = 0; $i -lt $newAgents.Length; $i++) {
    $row = 11 + $i
    $ws3.Cells.Item($row, 2).Value = "Agent"
    $ws3.Cells.Item($row, 3).Value = $newAgents[$i]
    $ws3.Cells.Item($row, 4).Value = $newProgress[$i]
}

# --- Selection/active-sheet bookkeeping --------------------------------
$ws3.Range("D11").Select()

$ws2.Activate()
$ws2.Range("C12").Select()
